# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos list.
# Numeric-looking Price values are entered with a leading apostrophe so
# Excel stores them as literal text (matching the sheet's existing
# inline-string cells) instead of silently re-parsing/rounding them as
# numbers (e.g. "12.50" -> 12.5, "0.100" -> 0.1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.622.71'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.669.05'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''599.23'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").Value = '''156.68'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.608'
$ws.Range("E8").Value = '  +3.38%  '
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("D10").Value = '''5.91'
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '''29.36'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '3.149.53'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '65.426.59'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '2.677.35'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '''12.50'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("D19").Value = '''4.81'
$ws.Range("E19").Value = '  -1.51%  '
$ws.Range("D20").Value = '''7.54'
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '''350.23'
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '''69.67'
$ws.Range("E24").Value = '  +3.49%  '
$ws.Range("D25").Value = '''9.73'
$ws.Range("E25").Value = '  +3.02%  '
$ws.Range("E26").Value = '  -3.23%  '
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").Value = '''1.60'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("D29").Value = '''8.11'
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").Value = '''541.53'
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("E32").Value = '  -2.97%  '
$ws.Range("D33").Value = '''1.76'
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("D34").Value = '''6.53'
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("D35").Value = '''5.46'
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("D36").Value = '''0.422'
$ws.Range("E36").Value = '  -2.57%  '
$ws.Range("D37").Value = '''20.42'
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = '''158.56'
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").Value = '''1.95'
$ws.Range("E40").Value = '  -2.55%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '''42.63'
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").Value = '''165.71'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("D45").Value = '''0.0611'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("E46").Value = '  -4.69%  '
$ws.Range("D47").Value = '''23.07'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").Value = '''0.644'
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("D50").Value = '''0.100'
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("D51").Value = '''19.93'
$ws.Range("E51").Value = '  +0.95%  '
